$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update per-unit bus voltage results (vm_pu) for the 380 kV case.
# Rows 2-25 correspond to bus indices 0-23; columns B:F and I:N hold
# the voltage magnitudes for the different buses (G/H unchanged).

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.062888517126526
$ws.Range("D2").Value = 1.064029046084612
$ws.Range("E2").Value = 1.066236055892689
$ws.Range("F2").Value = 1.076036832045149
$ws.Range("I2").Value = 1.048914487346738
$ws.Range("J2").Value = 1.067856720986985
$ws.Range("K2").Value = 1.066746087968209
$ws.Range("L2").Value = 1.06894714811469
$ws.Range("M2").Value = 1.07872182197649
$ws.Range("N2").Value = 1.069373200015664

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.064445966624521
$ws.Range("D3").Value = 1.065210333627218
$ws.Range("E3").Value = 1.067472314983088
$ws.Range("F3").Value = 1.077334726942754
$ws.Range("I3").Value = 1.049332773152549
$ws.Range("J3").Value = 1.069065889754938
$ws.Range("K3").Value = 1.067741583439468
$ws.Range("L3").Value = 1.069997910833187
$ws.Range("M3").Value = 1.079835969544343
$ws.Range("N3").Value = 1.070584085941962

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.065452247054846
$ws.Range("D4").Value = 1.065973249136241
$ws.Range("E4").Value = 1.068271595220965
$ws.Range("F4").Value = 1.078173562194179
$ws.Range("I4").Value = 1.0496011911611
$ws.Range("J4").Value = 1.069846406364954
$ws.Range("K4").Value = 1.068383705152092
$ws.Range("L4").Value = 1.070676597471806
$ws.Range("M4").Value = 1.08055535313931
$ws.Range("N4").Value = 1.07136571097509

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.065874937700264
$ws.Range("D5").Value = 1.066293635071169
$ws.Range("E5").Value = 1.068607458563858
$ws.Range("F5").Value = 1.078525976705842
$ws.Range("I5").Value = 1.049713499890243
$ws.Range("J5").Value = 1.070174087390779
$ws.Range("K5").Value = 1.068653171688098
$ws.Range("L5").Value = 1.07096162771285
$ws.Range("M5").Value = 1.080857416985498
$ws.Range("N5").Value = 1.071693857345557

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.065945889045089
$ws.Range("D6").Value = 1.066347409265251
$ws.Range("E6").Value = 1.068663842590583
$ws.Range("F6").Value = 1.078585135205271
$ws.Range("I6").Value = 1.049732325759813
$ws.Range("J6").Value = 1.070229080366193
$ws.Range("K6").Value = 1.068698388208524
$ws.Range("L6").Value = 1.0710094686811
$ws.Range("M6").Value = 1.080908113566002
$ws.Range("N6").Value = 1.071748928417305

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.065457896432034
$ws.Range("D7").Value = 1.065977531493362
$ws.Range("E7").Value = 1.068276083643293
$ws.Range("F7").Value = 1.078178272078794
$ws.Range("I7").Value = 1.049602693932097
$ws.Range("J7").Value = 1.069850786608587
$ws.Range("K7").Value = 1.068387307663097
$ws.Range("L7").Value = 1.070680407192722
$ws.Range("M7").Value = 1.08055939075966
$ws.Range("N7").Value = 1.071370097439171

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.063415178495811
$ws.Range("D8").Value = 1.064428572311171
$ws.Range("E8").Value = 1.066653995197583
$ws.Range("F8").Value = 1.076475669652132
$ws.Range("I8").Value = 1.049056314352629
$ws.Range("J8").Value = 1.068265761437063
$ws.Range("K8").Value = 1.067082943547432
$ws.Range("L8").Value = 1.069302514878864
$ws.Range("M8").Value = 1.079098675982162
$ws.Range("N8").Value = 1.069782821350099

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.059803854420089
$ws.Range("D9").Value = 1.061687730262989
$ws.Range("E9").Value = 1.063790398247757
$ws.Range("F9").Value = 1.073467659867943
$ws.Range("I9").Value = 1.048076258918646
$ws.Range("J9").Value = 1.065457956802723
$ws.Range("K9").Value = 1.064768727092218
$ws.Range("L9").Value = 1.066864907786669
$ws.Range("M9").Value = 1.076512669745449
$ws.Range("N9").Value = 1.066971029311119

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.057387876464994
$ws.Range("D10").Value = 1.059852538083807
$ws.Range("E10").Value = 1.061877509533308
$ws.Range("F10").Value = 1.071456749678247
$ws.Range("I10").Value = 1.047411136005879
$ws.Range("J10").Value = 1.063575774134568
$ws.Range("K10").Value = 1.063215024454911
$ws.Range("L10").Value = 1.065233129099982
$ws.Range("M10").Value = 1.074780271620991
$ws.Range("N10").Value = 1.065086173727686

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056339620589929
$ws.Range("D11").Value = 1.059055919934429
$ws.Range("E11").Value = 1.0610482373758
$ws.Range("F11").Value = 1.070584609855161
$ws.Range("I11").Value = 1.047120310165784
$ws.Range("J11").Value = 1.062758239065506
$ws.Range("K11").Value = 1.062539604279552
$ws.Range("L11").Value = 1.064524900882916
$ws.Range("M11").Value = 1.074028065964445
$ws.Range("N11").Value = 1.064267477665053

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.055949923261043
$ws.Range("D12").Value = 1.058759718654999
$ws.Range("E12").Value = 1.060740056202146
$ws.Range("F12").Value = 1.07026044127092
$ws.Range("I12").Value = 1.047011857448378
$ws.Range("J12").Value = 1.062454181649686
$ws.Range("K12").Value = 1.062288318200809
$ws.Range("L12").Value = 1.064261579223853
$ws.Range("M12").Value = 1.073748347041193
$ws.Range("N12").Value = 1.063962988452827

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.056033529679236
$ws.Range("D13").Value = 1.058823268593199
$ws.Range("E13").Value = 1.060806169103432
$ws.Range("F13").Value = 1.070329986455389
$ws.Range("I13").Value = 1.047035140296108
$ws.Range("J13").Value = 1.062519420694578
$ws.Range("K13").Value = 1.062342238336558
$ws.Range("L13").Value = 1.064318074198666
$ws.Range("M13").Value = 1.0738083621057
$ws.Range("N13").Value = 1.064028320144645

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056307414843236
$ws.Range("D14").Value = 1.059031442037466
$ws.Range("E14").Value = 1.061022766168385
$ws.Range("F14").Value = 1.07055781844279
$ws.Range("I14").Value = 1.047111354158877
$ws.Range("J14").Value = 1.062733113551904
$ws.Range("K14").Value = 1.06251884118862
$ws.Range("L14").Value = 1.064503139847423
$ws.Range("M14").Value = 1.074004950796596
$ws.Range("N14").Value = 1.064242316470339

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.056476120721569
$ws.Range("D15").Value = 1.059159664385817
$ws.Range("E15").Value = 1.061156198385439
$ws.Range("F15").Value = 1.070698164302158
$ws.Range("I15").Value = 1.047158255359667
$ws.Range("J15").Value = 1.062864725077387
$ws.Range("K15").Value = 1.062627598199448
$ws.Range("L15").Value = 1.064617131066608
$ws.Range("M15").Value = 1.074126033494285
$ws.Range("N15").Value = 1.064374114899286

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057457400232097
$ws.Range("D16").Value = 1.059905365020152
$ws.Range("E16").Value = 1.061932524503745
$ws.Range("F16").Value = 1.071514600639788
$ws.Range("I16").Value = 1.047430377430518
$ws.Range("J16").Value = 1.063629977163614
$ws.Range("K16").Value = 1.063259793398936
$ws.Range("L16").Value = 1.065280096516227
$ws.Range("M16").Value = 1.074830149007311
$ws.Range("N16").Value = 1.06514045373125

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.058072356142015
$ws.Range("D17").Value = 1.060372592133009
$ws.Range("E17").Value = 1.06241922744539
$ws.Range("F17").Value = 1.072026349664164
$ws.Range("I17").Value = 1.047600314522871
$ws.Range("J17").Value = 1.064109315102579
$ws.Range("K17").Value = 1.063655637384861
$ws.Range("L17").Value = 1.065695509883389
$ws.Range("M17").Value = 1.075271264844998
$ws.Range("N17").Value = 1.065620472385069

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.058430845487375
$ws.Range("D18").Value = 1.060644928357698
$ws.Range("E18").Value = 1.062703018980935
$ws.Range("F18").Value = 1.072324709327114
$ws.Range("I18").Value = 1.047699163750146
$ws.Range("J18").Value = 1.064388660468911
$ws.Range("K18").Value = 1.063886270549509
$ws.Range("L18").Value = 1.065937653921396
$ws.Range("M18").Value = 1.075528361327177
$ws.Range("N18").Value = 1.065900214453863

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.058553046774422
$ws.Range("D19").Value = 1.060737756025343
$ws.Range("E19").Value = 1.062799768773515
$ws.Range("F19").Value = 1.072426419573197
$ws.Range("I19").Value = 1.04773282269867
$ws.Range("J19").Value = 1.064483868793217
$ws.Range("K19").Value = 1.063964867262482
$ws.Range("L19").Value = 1.066020191873988
$ws.Range("M19").Value = 1.075615991019801
$ws.Range("N19").Value = 1.06599555798491

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.058006398367076
$ws.Range("D20").Value = 1.060322482725687
$ws.Range("E20").Value = 1.062367018594661
$ws.Range("F20").Value = 1.071971457834634
$ws.Range("I20").Value = 1.047582110057409
$ws.Range("J20").Value = 1.064057912004324
$ws.Range("K20").Value = 1.063613193586158
$ws.Range("L20").Value = 1.065650956501205
$ws.Range("M20").Value = 1.075223957908489
$ws.Range("N20").Value = 1.065568996288519

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056226771633162
$ws.Range("D21").Value = 1.058970148528619
$ws.Range("E21").Value = 1.06095898798762
$ws.Range("F21").Value = 1.070490733642948
$ws.Range("I21").Value = 1.047088922887941
$ws.Range("J21").Value = 1.062670197114336
$ws.Range("K21").Value = 1.062466847276018
$ws.Range("L21").Value = 1.064448649684493
$ws.Range("M21").Value = 1.073947069104504
$ws.Range("N21").Value = 1.064179310684214

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055105944661314
$ws.Range("D22").Value = 1.058118132951879
$ws.Range("E22").Value = 1.060072817762038
$ws.Range("F22").Value = 1.069558484784058
$ws.Range("I22").Value = 1.046776364300809
$ws.Range("J22").Value = 1.06179543436547
$ws.Range("K22").Value = 1.061743748006319
$ws.Range("L22").Value = 1.063691238553622
$ws.Range("M22").Value = 1.073142405758549
$ws.Range("N22").Value = 1.063303305671908

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.055700300073919
$ws.Range("D23").Value = 1.05856997043004
$ws.Range("E23").Value = 1.060542679156847
$ws.Range("F23").Value = 1.070052808984471
$ws.Range("I23").Value = 1.046942292833554
$ws.Range("J23").Value = 1.062259378484042
$ws.Range("K23").Value = 1.06212730092347
$ws.Range("L23").Value = 1.06409289789748
$ws.Range("M23").Value = 1.073569148714777
$ws.Range("N23").Value = 1.063767908644343

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.05803620247712
$ws.Range("D24").Value = 1.060345125595566
$ws.Range("E24").Value = 1.062390609817946
$ws.Range("F24").Value = 1.0719962615067
$ws.Range("I24").Value = 1.047590336712959
$ws.Range("J24").Value = 1.064081139607037
$ws.Range("K24").Value = 1.063632372902861
$ws.Range("L24").Value = 1.0656710887501
$ws.Range("M24").Value = 1.075245334491712
$ws.Range("N24").Value = 1.065592256877091

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.060738917065878
$ws.Range("D25").Value = 1.062397684177743
$ws.Range("E25").Value = 1.06453136000199
$ws.Range("F25").Value = 1.074246260588754
$ws.Range("I25").Value = 1.048331687529579
$ws.Range("J25").Value = 1.066185631997808
$ws.Range("K25").Value = 1.065368905450313
$ws.Range("L25").Value = 1.067496250458784
$ws.Range("M25").Value = 1.077182671789809
$ws.Range("N25").Value = 1.067699737888447

